# Update "想去人数" (number of people interested) values (column F)
# across the four worksheets, per the commit's regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1566
$ws1.Range("F5").Value = 227
$ws1.Range("F6").Value = 61
$ws1.Range("F7").Value = 1125
$ws1.Range("F8").Value = 724
$ws1.Range("F9").Value = 774
$ws1.Range("F10").Value = 1394
$ws1.Range("F12").Value = 1025
$ws1.Range("F13").Value = 30
$ws1.Range("F16").Value = 44
$ws1.Range("F17").Value = 439
$ws1.Range("F18").Value = 12
$ws1.Range("F21").Value = 544
$ws1.Range("F24").Value = 236
$ws1.Range("F25").Value = 171

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 7
$ws2.Range("F5").Value = 256
$ws2.Range("F7").Value = 140

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 212

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 212
$ws4.Range("F4").Value = 1566
$ws4.Range("F5").Value = 7
$ws4.Range("F7").Value = 227
$ws4.Range("F9").Value = 61
$ws4.Range("F10").Value = 1125
$ws4.Range("F11").Value = 724
$ws4.Range("F12").Value = 774
$ws4.Range("F13").Value = 1394
$ws4.Range("F15").Value = 1025
$ws4.Range("F16").Value = 30
$ws4.Range("F19").Value = 44
$ws4.Range("F20").Value = 439
$ws4.Range("F21").Value = 12
$ws4.Range("F23").Value = 256
$ws4.Range("F27").Value = 140
$ws4.Range("F28").Value = 140
$ws4.Range("F29").Value = 544
$ws4.Range("F32").Value = 236
$ws4.Range("F34").Value = 171
